$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.770.78'
$ws.Range("E2").Value = '  +2.73%  '
$ws.Range("D3").Value = '1.863.42'
$ws.Range("E3").Value = '  +2.31%  '
$ws.Range("D4").Value = '1.039'
$ws.Range("E4").Value = '  +2.74%  '
$ws.Range("D5").Value = '325.08'
$ws.Range("E5").Value = '  +3.52%  '
$ws.Range("D6").Value = '1.035'
$ws.Range("E7").Value = '  +2.23%  '
$ws.Range("D8").Value = '0.3794'
$ws.Range("E8").Value = '  +2.39%  '
$ws.Range("D9").Value = '0.07458'
$ws.Range("E9").Value = '  +2.41%  '
$ws.Range("D10").Value = '0.8839'
$ws.Range("E10").Value = '  +1.35%  '
$ws.Range("D11").Value = '21.76'
$ws.Range("E11").Value = '  +1.89%  '
$ws.Range("D12").Value = '1.881.76'
$ws.Range("E12").Value = '  -14.07%  '
$ws.Range("D13").Value = '5.555'
$ws.Range("E13").Value = '  +2.38%  '
$ws.Range("D14").Value = '6.748'
$ws.Range("E14").Value = '  +1.42%  '
$ws.Range("E15").Value = '  +3.30%  '
$ws.Range("D16").Value = '83.68'
$ws.Range("E16").Value = '  +3.02%  '
$ws.Range("E17").Value = '  +2.32%  '
$ws.Range("D18").Value = '0.000009135'
$ws.Range("E18").Value = '  +2.15%  '
$ws.Range("D19").Value = '1.034'
$ws.Range("E19").Value = '  +2.51%  '
$ws.Range("D20").Value = '15.53'
$ws.Range("E20").Value = '  +1.18%  '
$ws.Range("D21").Value = '27.801.73'
$ws.Range("E21").Value = '  +2.59%  '
$ws.Range("D22").Value = '5.316'
$ws.Range("E22").Value = '  +1.84%  '
$ws.Range("D23").Value = '11.38'
$ws.Range("E23").Value = '  +3.17%  '
$ws.Range("D24").Value = '1.965'
$ws.Range("E24").Value = '  +3.79%  '
$ws.Range("D25").Value = '158.33'
$ws.Range("E26").Value = '  +2.33%  '
$ws.Range("D27").Value = '1.989'
$ws.Range("E27").Value = '  +3.14%  '
$ws.Range("D28").Value = '5.317'
$ws.Range("E28").Value = '  +1.28%  '
$ws.Range("D29").Value = '117.54'
$ws.Range("E29").Value = '  +2.13%  '
$ws.Range("D30").Value = '0.09105'
$ws.Range("E30").Value = '  +1.33%  '
$ws.Range("D31").Value = '1.217'
$ws.Range("E31").Value = '  +3.42%  '
$ws.Range("E32").Value = '  +3.54%  '
$ws.Range("D33").Value = '3.076'
$ws.Range("E33").Value = '  +8.99%  '
$ws.Range("D34").Value = '4.582'
$ws.Range("E34").Value = '  +3.13%  '
$ws.Range("E35").Value = '  +2.58%  '
$ws.Range("D36").Value = '1.165'
$ws.Range("E36").Value = '  +3.44%  '
$ws.Range("D37").Value = '0.01996'
$ws.Range("E37").Value = '  +3.35%  '
$ws.Range("D38").Value = '0.05356'
$ws.Range("E38").Value = '  +2.00%  '
$ws.Range("D39").Value = '0.5203'
$ws.Range("E39").Value = '  +1.09%  '
$ws.Range("D40").Value = '2.843'
$ws.Range("E40").Value = '  +3.46%  '
$ws.Range("E41").Value = '  +2.40%  '
$ws.Range("D42").Value = '6.901'
$ws.Range("E42").Value = '  +5.89%  '
$ws.Range("D43").Value = '8.716'
$ws.Range("E43").Value = '  +4.27%  '
$ws.Range("D44").Value = '10.73'
$ws.Range("E44").Value = '  +2.97%  '
$ws.Range("D45").Value = '109.62'
$ws.Range("E45").Value = '  +1.96%  '
$ws.Range("D46").Value = '1.721'
$ws.Range("E46").Value = '  +3.78%  '
$ws.Range("D47").Value = '0.4713'
$ws.Range("E47").Value = '  +2.40%  '
$ws.Range("D48").Value = '0.06438'
$ws.Range("E48").Value = '  +1.79%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.880'
$ws.Range("E49").Value = '  +3.08%  '
$ws.Range("D50").Value = '39.97'
$ws.Range("E50").Value = '  +4.52%  '
$ws.Range("D51").Value = '64.56'
$ws.Range("E51").Value = '  +0.92%  '
